# Generate Report for Handoff
# Updates the "b.md" row (row 3) across the Overview, zh-cn and de-de sheets
# to reflect that a new handoff package has been generated.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: b.md row (row 3)
#   E3 (zh-cn status) / F3 (de-de status) -> "Ready for handoff"
#   G3 (Latest HO Xliff Generate Date)    -> "2016-09-05 06:43:09"
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-05 06:43:09"

# ---------------------------------------------------------------------
# zh-cn sheet: b.md row (row 3)
#   C3 (Status)                  -> "Ready for handoff"
#   F3 (Content Duplicate)       -> "False"
#   G3 (Latest Handoff File)     -> new zh-cn xlf handoff file name
#   H3 (Latest Handoff Datetime) -> "2016-09-05 06:42:59"
#   P3 (Error Detail)            -> version-mismatch message
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-09-05 06:42:59"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c6b94f2e15aca95e9f021b41c6043dff12c2fce4/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/23248a9eb78b05a954ead04d28d0fe22445f06ad/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 40

# ---------------------------------------------------------------------
# de-de sheet: b.md row (row 3)
#   C3 (Status)                  -> "Ready for handoff"
#   F3 (Content Duplicate)       -> "False"
#   G3 (Latest Handoff File)     -> new de-de xlf handoff file name
#   H3 (Latest Handoff Datetime) -> "2016-09-05 06:43:09"
#   P3 (Error Detail)            -> version-mismatch message
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-09-05 06:43:09"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c6b94f2e15aca95e9f021b41c6043dff12c2fce4/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/23248a9eb78b05a954ead04d28d0fe22445f06ad/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 40
